$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the data range so numeric-looking strings
# (e.g. "0.991", "1.85") are stored as text, matching the source data
# which uses inline strings for all of columns B:E.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '60.723.53'
$ws.Range('E2').Value = '  -3.66%  '
$ws.Range('D3').Value = '2.509.77'
$ws.Range('E3').Value = '  -6.29%  '
$ws.Range('E4').Value = '  +0.43%  '
$ws.Range('D5').Value = '536.89'
$ws.Range('E5').Value = '  -2.36%  '
$ws.Range('D6').Value = '148.71'
$ws.Range('E6').Value = '  -5.75%  '
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('D8').Value = '0.571'
$ws.Range('E8').Value = '  -2.46%  '
$ws.Range('D9').Value = '0.100'
$ws.Range('E9').Value = '  -4.59%  '
$ws.Range('D10').Value = '0.160'
$ws.Range('E10').Value = '  -1.26%  '
$ws.Range('E11').Value = '  +5.70%  '
$ws.Range('D12').Value = '0.354'
$ws.Range('E12').Value = '  -3.36%  '
$ws.Range('D13').Value = '2.992.33'
$ws.Range('E13').Value = '  -4.99%  '
$ws.Range('D14').Value = '24.72'
$ws.Range('E14').Value = '  -5.13%  '
$ws.Range('D15').Value = '60.633.33'
$ws.Range('E15').Value = '  -3.58%  '
$ws.Range('D16').Value = '0.0000139'
$ws.Range('E16').Value = '  -3.54%  '
$ws.Range('D17').Value = '2.544.85'
$ws.Range('E17').Value = '  -5.02%  '
$ws.Range('D18').Value = '11.30'
$ws.Range('E18').Value = '  -5.17%  '
$ws.Range('D19').Value = '4.40'
$ws.Range('E19').Value = '  -3.46%  '
$ws.Range('D20').Value = '329.14'
$ws.Range('E20').Value = '  -4.01%  '
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('D22').Value = '5.85'
$ws.Range('E22').Value = '  -7.11%  '
$ws.Range('D23').Value = '0.476'
$ws.Range('E23').Value = '  -5.66%  '
$ws.Range('D24').Value = '61.88'
$ws.Range('E24').Value = '  -2.30%  '
$ws.Range('D25').Value = '0.163'
$ws.Range('E25').Value = '  -2.54%  '
$ws.Range('D26').Value = '0.991'
$ws.Range('E26').Value = '  -0.92%  '
$ws.Range('D27').Value = '7.88'
$ws.Range('E27').Value = '  -2.95%  '
$ws.Range('D28').Value = '7.03'
$ws.Range('E28').Value = '  +0.43%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').Value = '1.85'
$ws.Range('E29').Value = '  -3.50%  '
$ws.Range('D30').Value = '1.29'
$ws.Range('E30').Value = '  -3.60%  '
$ws.Range('B31').Value = 'PEPE'
$ws.Range('C31').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D31').Value = '0.0₃0787'
$ws.Range('E31').Value = '  -7.77%  '
$ws.Range('D32').Value = '161.08'
$ws.Range('E32').Value = '  -3.07%  '
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('D34').Value = '4.64'
$ws.Range('E34').Value = '  -3.39%  '
$ws.Range('D35').Value = '18.65'
$ws.Range('E35').Value = '  -4.28%  '
$ws.Range('D36').Value = '1.36'
$ws.Range('E36').Value = '  -4.73%  '
$ws.Range('D37').Value = '1.72'
$ws.Range('E37').Value = '  -2.95%  '
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').Value = '316.51'
$ws.Range('E38').Value = '  -6.65%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = '5.80'
$ws.Range('E39').Value = '  -5.09%  '
$ws.Range('D40').Value = '0.867'
$ws.Range('E40').Value = '  -7.04%  '
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').Value = '36.99'
$ws.Range('E41').Value = '  -2.78%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').Value = '3.81'
$ws.Range('E42').Value = '  -2.86%  '
$ws.Range('D43').Value = '1.01'
$ws.Range('E43').Value = '  +0.97%  '
$ws.Range('B44').Value = 'WhiteBITCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D44').Value = '10.89'
$ws.Range('E44').Value = '  -1.38%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = '19.91'
$ws.Range('E45').Value = '  -3.71%  '
$ws.Range('D46').Value = '0.595'
$ws.Range('E46').Value = '  -3.14%  '
$ws.Range('D47').Value = '0.0947'
$ws.Range('E47').Value = '  -2.31%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '18.98'
$ws.Range('E48').Value = '  -6.60%  '
$ws.Range('B49').Value = 'Hedera'
$ws.Range('C49').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D49').Value = '0.0526'
$ws.Range('E49').Value = '  -5.99%  '
$ws.Range('D50').Value = '0.0230'
$ws.Range('E50').Value = '  -3.93%  '
$ws.Range('D51').Value = '1.997.45'
$ws.Range('E51').Value = '  -3.73%  '

# Restore default (Normal) style so we don't leave stray number formatting
# behind on cells that did not originally carry an explicit style.
$dataRange.Style = "Normal"
